# "Wrote name, email and repo link 2" -- adds a "fork link" (label in D2,
# URL+hyperlink in C3) and a "project link" label (D3) to the table that
# already has name / email / repo link in row 1-2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label/value cells.
$ws.Range("D2").Value = "fork link"
$ws.Range("C3").Value = "https://github.com/THOMASSAAD/Open_Source"
$ws.Range("D3").Value = "project link"

# Turn C3 into a live hyperlink, same as the existing email/repo links.
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/THOMASSAAD/Open_Source")

# Match C3's visual style to the other hyperlink cell in column C (C2) so
# it reuses the same cell style instead of a freshly-minted one.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen column C and size the new column D to fit the new content.
$ws.Columns.Item(3).ColumnWidth = 46.833333333333336
$ws.Columns.Item(4).ColumnWidth = 13.666666666666666

# Matches the saved selection/active cell recorded in the workbook.
$ws.Range("F9").Select()
